# Update the cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even if it looks numeric
# (e.g. "211.98"), without leaving a residual NumberFormat/Style on
# the cell - mirrors the original file's plain (un-styled) text cells.
function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Column D = Price, Column E = Volume(1h)
# Each entry: row number, new D value (or $null if unchanged), new E value (or $null if unchanged)
$updates = @(
    @{Row = 2;  D = "27.915.89";  E = $null},
    @{Row = 3;  D = "1.636.93";   E = "  +0.25%  "},
    @{Row = 4;  D = $null;        E = "  +0.00%  "},
    @{Row = 5;  D = "211.98";     E = "  +0.10%  "},
    @{Row = 6;  D = $null;        E = "  -0.45%  "},
    @{Row = 7;  D = $null;        E = "  -0.01%  "},
    @{Row = 8;  D = "23.48";      E = "  +0.95%  "},
    @{Row = 9;  D = $null;        E = "  -0.32%  "},
    @{Row = 10; D = $null;        E = "  -0.27%  "},
    @{Row = 11; D = $null;        E = "  +0.38%  "},
    @{Row = 12; D = "1.868.45";   E = "  +0.22%  "},
    @{Row = 13; D = "1.638.11";   E = "  +0.37%  "},
    @{Row = 14; D = $null;        E = "  -0.87%  "},
    @{Row = 15; D = "0.565";      E = "  -0.12%  "},
    @{Row = 16; D = "65.39";      E = "  +0.26%  "},
    @{Row = 17; D = "27.926.66";  E = "  +0.03%  "},
    @{Row = 18; D = "229.73";     E = "  -0.32%  "},
    @{Row = 19; D = "7.88";       E = "  +4.91%  "},
    @{Row = 20; D = $null;        E = "  -0.24%  "},
    @{Row = 21; D = "1.00";       E = "  -0.02%  "},
    @{Row = 23; D = "10.21";      E = "  -1.49%  "},
    @{Row = 24; D = $null;        E = "  +0.60%  "},
    @{Row = 25; D = "156.27";     E = "  +1.17%  "},
    @{Row = 26; D = "6.98";       E = "  +0.29%  "},
    @{Row = 27; D = $null;        E = "  +0.04%  "},
    @{Row = 28; D = "15.58";      E = "  -0.42%  "},
    @{Row = 29; D = $null;        E = "  +0.11%  "},
    @{Row = 30; D = "1.19";       E = "  +0.27%  "},
    @{Row = 31; D = $null;        E = "  -0.23%  "},
    @{Row = 32; D = $null;        E = "  +1.06%  "},
    @{Row = 33; D = $null;        E = "  +1.38%  "},
    @{Row = 34; D = "1.404.59";   E = "  +0.26%  "},
    @{Row = 35; D = $null;        E = "  +3.18%  "},
    @{Row = 36; D = $null;        E = "  +0.62%  "},
    @{Row = 37; D = $null;        E = "  -0.77%  "},
    @{Row = 38; D = $null;        E = "  +0.71%  "},
    @{Row = 39; D = "0.559";      E = "  -0.23%  "},
    @{Row = 40; D = "0.855";      E = "  -1.87%  "},
    @{Row = 41; D = "1.00";       E = "  +0.05%  "},
    @{Row = 42; D = $null;        E = "  -1.24%  "},
    @{Row = 43; D = $null;        E = "  +2.32%  "},
    @{Row = 44; D = "66.28";      E = "  -0.77%  "},
    @{Row = 45; D = "5.48";       E = "  -1.19%  "},
    @{Row = 46; D = "1.777.36";   E = "  +0.13%  "},
    @{Row = 47; D = $null;        E = "  -2.69%  "},
    @{Row = 48; D = "88.82";      E = "  +1.31%  "},
    @{Row = 49; D = "0.103";      E = "  +2.46%  "},
    @{Row = 50; D = "0.0505";     E = "  -0.28%  "},
    @{Row = 51; D = "7.63";       E = "  +2.19%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextValue $r 4 $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
